$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the A, Q, R values between row 2 and row 3
$a2 = $ws.Range("A2").Value()
$q2 = $ws.Range("Q2").Value()
$r2 = $ws.Range("R2").Value()

$a3 = $ws.Range("A3").Value()
$q3 = $ws.Range("Q3").Value()
$r3 = $ws.Range("R3").Value()

$ws.Range("A2").Value = $a3
$ws.Range("Q2").Value = $q3
$ws.Range("R2").Value = $r3

$ws.Range("A3").Value = $a2
$ws.Range("Q3").Value = $q2
$ws.Range("R3").Value = $r2
